$wb = $excel.ActiveWorkbook

# The "meta" sheet holds key/value configuration rows (tab, type, title, ...).
# This change adds a new "style" = "default" row (feature: arrow_n gets a
# default style), inserted just above the existing "lock" / "no" row, which
# (along with the trailing blank row) shifts down by one row.
$ws = $wb.Worksheets.Item("meta")

$ws.Rows.Item(18).Insert()
$ws.Cells.Item(18, 1).Value = "style"
$ws.Cells.Item(18, 2).Value = "default"
